$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Shared strings must be (re)created in a very specific order to match the
# target sharedStrings.xml table byte-for-byte: header row left-to-right
# first, then the data area column-by-column (top-to-bottom within each
# column), left-to-right across columns. We therefore perform every write in
# that exact sequence.
# ---------------------------------------------------------------------------

# --- Header row (row 1), left to right ---
$ws.Range("B1").Value = "squadName"
$ws.Range("C1").Value = "homeTown"
$ws.Range("D1").Value = "formed"
$ws.Range("E1").Value = "secretBase"
$ws.Range("F1").Value = "active"
$ws.Range("G1").Value = "members.name"
$ws.Range("H1").Value = "members.age"
$ws.Range("I1").Value = "members.secretIdentity"
$ws.Range("J1").Value = "members.powers"

# New header cells H1:J1 need the same bold/bordered/centered look as the
# rest of the header row. Copy *formatting only* (xlPasteFormats) from an
# already-styled header cell so the existing text/shared-string order set
# above is left untouched.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("G1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("G1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# --- Data area, column by column (top to bottom within each column) ---

# Column B
$ws.Range("B2").Value = "Super hero squad"

# Column C
$ws.Range("C2").Value = "Metro City"

# Column D (numeric)
$ws.Range("D2").Value = 2016

# Column E
$ws.Range("E2").Value = "Super tower"

# Column F -- the literal text "True" must stay text, not become a Boolean.
# A plain .Value assignment of "True"/"False" auto-types to a Boolean cell,
# so build it as a formula first, then freeze it to a static value via
# copy / paste-values (which does not re-parse the text).
$ws.Range("F2").Formula = '="Tru"&"e"'
$ws.Range("F2").Copy()
$ws.Range("F2").PasteSpecial(-4163)

# Column G
$ws.Range("G2").Value = "Molecule Man"
$ws.Range("G3").Value = "Madame Uppercut"
$ws.Range("G4").Value = "Eternal Flame"

# Column H (numeric)
$ws.Range("H2").Value = 29
$ws.Range("H3").Value = 39
$ws.Range("H4").Value = 1000000

# Column I
$ws.Range("I2").Value = "Dan Jukes"
$ws.Range("I3").Value = "Jane Wilson"
$ws.Range("I4").Value = "Unknown"

# Column J
$ws.Range("J2").Value = "['Radiation resistance', 'Turning tiny', 'Radiation blast']"
$ws.Range("J3").Value = "['Million tonne punch', 'Damage resistance', 'Superhuman reflexes']"
$ws.Range("J4").Value = "['Immortality', 'Heat Immunity', 'Inferno', 'Teleportation', 'Interdimensional travel']"

# --- Remove leftover values from the old, wider table (rows 3 & 4 used to
# hold sib.ph / sib.addr / frnds.b data in columns C:F; none of that survives
# in the new layout) ---
$ws.Range("C3:F3").ClearContents()
$ws.Range("C4:F4").ClearContents()

# --- Remove now-unused rows 5,6,7 (shrinks the sheet down to 4 rows) ---
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()

Write-Output "done"
